$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.806.67'
$ws.Range('E2').Value = '  +6.59%  '
$ws.Range('D3').Value = '3.354.22'
$ws.Range('E3').Value = '  +3.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '413.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '111.71'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.588'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.641'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.65'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0996'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.144'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.26%  '
$ws.Range('D13').Value = '3.864.91'
$ws.Range('E13').Value = '  +2.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.44'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.89'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.40%  '
$ws.Range('D16').Value = '3.358.60'
$ws.Range('E16').Value = '  +3.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.05'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.48%  '
$ws.Range('D18').Value = '60.267.03'
$ws.Range('E18').Value = '  +6.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.73'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.38'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000111'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.13'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '302.64'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.19'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.86'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '28.76'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.83%  '
$ws.Range('E28').Value = '  +2.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.05'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.180'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.05%  '
$ws.Range('E31').Value = '  +25.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.116'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.59%  '
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '39.65'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0508'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.42'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.14'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.997'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.40'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '137.79'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.61%  '
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.298'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.90%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.124'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.92'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.96'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.47%  '
$ws.Range('E47').Value = '  +8.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.41'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.76%  '
$ws.Range('D49').Value = '2.181.37'
$ws.Range('E49').Value = '  +1.99%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.43'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.98%  '
$ws.Range('E51').Value = '  -1.36%  '
